$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 666.5172
$ws.Range("I33").Value = 353.875
$ws.Range("J33").Value = 2167.2
$ws.Range("K33").Value = 353.875
$ws.Range("L33").Value = 2167.2
$ws.Range("M33").Value = -124.875
$ws.Range("N33").Value = -2625.2

$ws.Range("H53").Value = 1158.6666
$ws.Range("I53").Value = 1678.5714
$ws.Range("J53").Value = 703.75
$ws.Range("K53").Value = 1678.5714
$ws.Range("L53").Value = 703.75
$ws.Range("M53").Value = -1041.5714
$ws.Range("N53").Value = -1977.75

$ws.Range("H62").Value = 2487.7144
$ws.Range("I62").Value = 2487.7144
$ws.Range("K62").Value = 2487.7144
$ws.Range("M62").Value = -1863.7144

$ws.Range("H65").Value = 2487.7144
$ws.Range("I65").Value = 2487.7144
$ws.Range("K65").Value = 12438.572
$ws.Range("M65").Value = -9318.572

$ws.Range("H112").Value = 1568.8334
$ws.Range("I112").Value = 775
$ws.Range("J112").Value = 1652.3948
$ws.Range("K112").Value = 2325
$ws.Range("L112").Value = 4957.1844
$ws.Range("M112").Value = -1217
$ws.Range("N112").Value = -7173.1844

$ws.Range("H115").Value = 1854.3572
$ws.Range("I115").Value = 250.125
$ws.Range("J115").Value = 3993.3333
$ws.Range("K115").Value = 750.375
$ws.Range("L115").Value = 11979.9999
$ws.Range("M115").Value = 816.625
$ws.Range("N115").Value = -15113.9999

$ws.Range("H129").Value = 1523.3906
$ws.Range("I129").Value = 492.42856
$ws.Range("J129").Value = 1650
$ws.Range("K129").Value = 1477.28568
$ws.Range("L129").Value = 4950
$ws.Range("M129").Value = 3522.71432
$ws.Range("N129").Value = -14950

$ws.Range("H137").Value = 994.5
$ws.Range("I137").Value = 1061.3846
$ws.Range("J137").Value = 915.4545000000001
$ws.Range("K137").Value = 3184.1538
$ws.Range("L137").Value = 2746.3635
$ws.Range("M137").Value = -634.1538
$ws.Range("N137").Value = -7846.3635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4630.56
$ws.Range("I32").Value = 4176.2637
$ws.Range("J32").Value = 9224
$ws.Range("K32").Value = 4176.2637
$ws.Range("L32").Value = 9224
$ws.Range("M32").Value = -3889.2637
$ws.Range("N32").Value = -9798

$ws.Range("H34").Value = 16023.333
$ws.Range("J34").Value = 16023.333
$ws.Range("L34").Value = 16023.333
$ws.Range("N34").Value = -16565.333

$ws.Range("H61").Value = 2103.6738
$ws.Range("I61").Value = 2156.875
$ws.Range("J61").Value = 1749
$ws.Range("K61").Value = 2156.875
$ws.Range("L61").Value = 1749
$ws.Range("M61").Value = -1944.875
$ws.Range("N61").Value = -2173

$ws.Range("H74").Value = 1514.8485
$ws.Range("I74").Value = 1203.9131
$ws.Range("K74").Value = 1203.9131
$ws.Range("M74").Value = -329.9131

$ws.Range("H77").Value = 1514.8485
$ws.Range("I77").Value = 1203.9131
$ws.Range("K77").Value = 6019.5655
$ws.Range("M77").Value = -1651.5655

$ws.Range("H122").Value = 5520.242
$ws.Range("I122").Value = 6464.154
$ws.Range("J122").Value = 2014.2858
$ws.Range("K122").Value = 19392.462
$ws.Range("L122").Value = 6042.857400000001
$ws.Range("M122").Value = -16942.462
$ws.Range("N122").Value = -10942.8574

$ws.Range("H132").Value = 5436530.5
$ws.Range("I132").Value = 7813869.5
$ws.Range("K132").Value = 23441608.5
$ws.Range("M132").Value = -23439078.5

$ws.Range("H136").Value = 2103.6738
$ws.Range("I136").Value = 2156.875
$ws.Range("J136").Value = 1749
$ws.Range("K136").Value = 6470.625
$ws.Range("L136").Value = 5247
$ws.Range("M136").Value = -3920.625
$ws.Range("N136").Value = -10347

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 20734.143
$ws.Range("I81").Value = 5000
$ws.Range("J81").Value = 23356.5
$ws.Range("K81").Value = 5000
$ws.Range("L81").Value = 23356.5
$ws.Range("M81").Value = -3939
$ws.Range("N81").Value = -25478.5

$ws.Range("H84").Value = 20734.143
$ws.Range("I84").Value = 5000
$ws.Range("J84").Value = 23356.5
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 70069.5
$ws.Range("M84").Value = -9696
$ws.Range("N84").Value = -80677.5

$ws.Range("H99").Value = 780.6087
$ws.Range("I99").Value = 737.9
$ws.Range("J99").Value = 813.46155
$ws.Range("K99").Value = 737.9
$ws.Range("L99").Value = 813.46155
$ws.Range("M99").Value = 760.1
$ws.Range("N99").Value = -3809.46155

$ws.Range("H105").Value = 3686.923
$ws.Range("I105").Value = 1676.25
$ws.Range("J105").Value = 4205.8066
$ws.Range("K105").Value = 1676.25
$ws.Range("L105").Value = 4205.8066
$ws.Range("M105").Value = 70.75
$ws.Range("N105").Value = -7699.8066

$ws.Range("H134").Value = 3926.52
$ws.Range("I134").Value = 2803.9285
$ws.Range("J134").Value = 5355.273
$ws.Range("K134").Value = 8411.7855
$ws.Range("L134").Value = 16065.819
$ws.Range("M134").Value = -5876.7855
$ws.Range("N134").Value = -21135.819

$ws.Range("H135").Value = 43874.75
$ws.Range("J135").Value = 43874.75
$ws.Range("L135").Value = 43874.75
$ws.Range("N135").Value = -54014.75

$ws.Range("H138").Value = 49381.125
$ws.Range("J138").Value = 49381.125
$ws.Range("L138").Value = 49381.125
$ws.Range("N138").Value = -59661.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 6400
$ws.Range("I36").Value = 6400
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 6400
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -6012
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 6400
$ws.Range("I40").Value = 6400
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6400
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6240
$ws.Range("N40").ClearContents()

$ws.Range("H105").Value = 2107.2307
$ws.Range("I105").Value = 2049.6
$ws.Range("J105").Value = 2299.3333
$ws.Range("K105").Value = 2049.6
$ws.Range("L105").Value = 2299.3333
$ws.Range("M105").Value = -302.5999999999999
$ws.Range("N105").Value = -5793.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 11007.25
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 11007.25
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 33021.75
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -38263.75

$ws.Range("H113").Value = 27778474
$ws.Range("I113").Value = 76923600
$ws.Range("J113").Value = 793.8261
$ws.Range("K113").Value = 230770800
$ws.Range("L113").Value = 2381.4783
$ws.Range("M113").Value = -230768630
$ws.Range("N113").Value = -6721.4783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28500
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 4000
$ws.Range("N70").Value = -4540

$ws.Range("H73").Value = 28500
$ws.Range("J73").Value = 4000
$ws.Range("L73").Value = 4000
$ws.Range("N73").Value = -5872

$ws.Range("H107").Value = 1115.1818
$ws.Range("I107").Value = 1424.5
$ws.Range("J107").Value = 744
$ws.Range("K107").Value = 1424.5
$ws.Range("L107").Value = 744
$ws.Range("M107").Value = 495.5
$ws.Range("N107").Value = -4584

$ws.Range("H113").Value = 112607.445
$ws.Range("I113").Value = 168127.83
$ws.Range("J113").Value = 1566.6666
$ws.Range("K113").Value = 168127.83
$ws.Range("L113").Value = 1566.6666
$ws.Range("M113").Value = -165957.83
$ws.Range("N113").Value = -5906.6666

$ws.Range("H122").Value = 2572.8462
$ws.Range("I122").Value = 2383.111
$ws.Range("J122").Value = 2999.75
$ws.Range("K122").Value = 7149.333
$ws.Range("L122").Value = 8999.25
$ws.Range("M122").Value = -4699.333
$ws.Range("N122").Value = -13899.25

$ws.Range("H123").Value = 22790.3
$ws.Range("J123").Value = 22790.3
$ws.Range("L123").Value = 22790.3
$ws.Range("N123").Value = -27690.3

$ws.Range("H132").Value = 7239.36
$ws.Range("I132").Value = 9016.117
$ws.Range("K132").Value = 27048.351
$ws.Range("M132").Value = -24518.351

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6950594.5
$ws.Range("I132").Value = 4305.3794
$ws.Range("J132").Value = 35728080
$ws.Range("K132").Value = 12916.1382
$ws.Range("L132").Value = 107184240
$ws.Range("M132").Value = -10386.1382
$ws.Range("N132").Value = -107189300

$ws.Range("H136").Value = 5469.7646
$ws.Range("I136").Value = 2319.0344
$ws.Range("J136").Value = 23744
$ws.Range("K136").Value = 6957.1032
$ws.Range("L136").Value = 71232
$ws.Range("M136").Value = -4407.1032
$ws.Range("N136").Value = -76332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1329.4265
$ws.Range("I132").Value = 1152.5186
$ws.Range("K132").Value = 3457.5558
$ws.Range("M132").Value = -927.5558000000001

$ws.Range("H133").Value = 37132.855
$ws.Range("J133").Value = 37132.855
$ws.Range("L133").Value = 37132.855
$ws.Range("N133").Value = -47252.855

$ws.Range("H136").Value = 1325.66
$ws.Range("I136").Value = 847.30554
$ws.Range("J136").Value = 2555.7144
$ws.Range("K136").Value = 2541.91662
$ws.Range("L136").Value = 7667.1432
$ws.Range("M136").Value = 8.083380000000034
$ws.Range("N136").Value = -12767.1432
